$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 206.46153
$ws.Range("I33").Value = 201
$ws.Range("J33").Value = 224.66667
$ws.Range("K33").Value = 201
$ws.Range("L33").Value = 224.66667
$ws.Range("M33").Value = 28
$ws.Range("N33").Value = -682.6666700000001

$ws.Range("H80").Value = 104302.35
$ws.Range("J80").Value = 5550.273
$ws.Range("L80").Value = 16650.819
$ws.Range("N80").Value = -18646.819

$ws.Range("H83").Value = 104302.35
$ws.Range("J83").Value = 5550.273
$ws.Range("L83").Value = 49952.457
$ws.Range("N83").Value = -59936.457

$ws.Range("H132").Value = 5661.8696
$ws.Range("I132").Value = 5011.25
$ws.Range("K132").Value = 15033.75
$ws.Range("M132").Value = -12503.75

$ws.Range("H135").Value = 5666.933
$ws.Range("I135").Value = 6309.909
$ws.Range("K135").Value = 56789.181
$ws.Range("M135").Value = -54254.181

$ws.Range("H139").Value = 158201.17
$ws.Range("J139").Value = 158201.17
$ws.Range("L139").Value = 158201.17
$ws.Range("N139").Value = -168481.17

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5223.839
$ws.Range("I61").Value = 5377
$ws.Range("K61").Value = 5377
$ws.Range("M61").Value = -5165

$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()

$ws.Range("H112").Value = 34450
$ws.Range("J112").Value = 34450
$ws.Range("L112").Value = 34450
$ws.Range("N112").Value = -37404

$ws.Range("H122").Value = 429354.34
$ws.Range("I122").Value = 3909.4285
$ws.Range("K122").Value = 11728.2855
$ws.Range("M122").Value = -9278.2855

$ws.Range("H136").Value = 5223.839
$ws.Range("I136").Value = 5377
$ws.Range("K136").Value = 16131
$ws.Range("M136").Value = -13581

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 1143.5
$ws.Range("J80").Value = 1159.4445
$ws.Range("L80").Value = 1159.4445
$ws.Range("N80").Value = -3155.4445

$ws.Range("H81").Value = 27180
$ws.Range("J81").Value = 27725
$ws.Range("L81").Value = 27725
$ws.Range("N81").Value = -29847

$ws.Range("H83").Value = 1143.5
$ws.Range("J83").Value = 1159.4445
$ws.Range("L83").Value = 5797.2225
$ws.Range("N83").Value = -15781.2225

$ws.Range("H84").Value = 27180
$ws.Range("J84").Value = 27725
$ws.Range("L84").Value = 83175
$ws.Range("N84").Value = -93783

$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

$ws.Range("H126").Value = 58832.668
$ws.Range("J126").Value = 59999.4
$ws.Range("L126").Value = 59999.4
$ws.Range("N126").Value = -69879.39999999999

$ws.Range("H135").Value = 90469.336
$ws.Range("J135").Value = 90469.336
$ws.Range("L135").Value = 90469.336
$ws.Range("N135").Value = -100609.336

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4016.8
$ws.Range("I31").Value = 1233.6
$ws.Range("K31").Value = 1233.6
$ws.Range("M31").Value = -938.5999999999999

$ws.Range("H34").Value = 4016.8
$ws.Range("I34").Value = 1233.6
$ws.Range("K34").Value = 1233.6
$ws.Range("M34").Value = -1031.6

$ws.Range("H86").Value = 12710.96
$ws.Range("I86").Value = 12265.056
$ws.Range("K86").Value = 12265.056
$ws.Range("M86").Value = -11142.056

$ws.Range("H89").Value = 12710.96
$ws.Range("I89").Value = 12265.056
$ws.Range("K89").Value = 61325.28
$ws.Range("M89").Value = -55709.28

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 177.52942
$ws.Range("I12").Value = 418
$ws.Range("J12").Value = 126
$ws.Range("K12").Value = 1254
$ws.Range("L12").Value = 378
$ws.Range("M12").Value = -1081
$ws.Range("N12").Value = -724

$ws.Range("H23").Value = 525
$ws.Range("I23").Value = 402.4
$ws.Range("K23").Value = 1207.2
$ws.Range("M23").Value = -972.1999999999998

$ws.Range("H131").Value = 6278.933
$ws.Range("J131").Value = 2047.2858
$ws.Range("L131").Value = 6141.857400000001
$ws.Range("N131").Value = -16221.8574

$ws.Range("H132").Value = 51309.9
$ws.Range("I132").Value = 728.4286
$ws.Range("J132").Value = 169333.33
$ws.Range("K132").Value = 6555.8574
$ws.Range("L132").Value = 1523999.97
$ws.Range("M132").Value = -4025.8574
$ws.Range("N132").Value = -1529059.97

$ws.Range("H140").Value = 11687.474
$ws.Range("I140").Value = 13838.2
$ws.Range("J140").Value = 3622.25
$ws.Range("K140").Value = 41514.60000000001
$ws.Range("L140").Value = 10866.75
$ws.Range("M140").Value = -36334.60000000001
$ws.Range("N140").Value = -21226.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 3930000
$ws.Range("J11").Value = 3288571.2
$ws.Range("L11").Value = 3288571.2
$ws.Range("N11").Value = -3288849.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 17971.285
$ws.Range("I13").Value = 20633.166
$ws.Range("J13").Value = 2000
$ws.Range("K13").Value = 20633.166
$ws.Range("L13").Value = 2000
$ws.Range("M13").Value = -20493.166
$ws.Range("N13").Value = -2280

$ws.Range("H46").Value = 4126.5713
$ws.Range("I46").Value = 965.6667
$ws.Range("J46").Value = 4988.636
$ws.Range("K46").Value = 965.6667
$ws.Range("L46").Value = 4988.636
$ws.Range("M46").Value = -777.6667
$ws.Range("N46").Value = -5364.636

$ws.Range("H68").Value = 5075
$ws.Range("I68").Value = 1933.3334
$ws.Range("K68").Value = 1933.3334
$ws.Range("M68").Value = -1184.3334

$ws.Range("H71").Value = 5075
$ws.Range("I71").Value = 1933.3334
$ws.Range("K71").Value = 9666.666999999999
$ws.Range("M71").Value = -5922.666999999999

$ws.Range("H132").Value = 384982.84
$ws.Range("I132").Value = 553791.5
$ws.Range("K132").Value = 1661374.5
$ws.Range("M132").Value = -1658844.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H6").Value = 1901.25
$ws.Range("I6").Value = 1005
$ws.Range("J6").Value = 2200
$ws.Range("K6").Value = 1005
$ws.Range("L6").Value = 2200
$ws.Range("M6").Value = -890
$ws.Range("N6").Value = -2430

$ws.Range("H62").Value = 93466.734
$ws.Range("I62").Value = 144789.58
$ws.Range("J62").Value = 4818.1816
$ws.Range("K62").Value = 144789.58
$ws.Range("L62").Value = 4818.1816
$ws.Range("M62").Value = -144165.58
$ws.Range("N62").Value = -6066.1816

$ws.Range("H65").Value = 93466.734
$ws.Range("I65").Value = 144789.58
$ws.Range("J65").Value = 4818.1816
$ws.Range("K65").Value = 723947.8999999999
$ws.Range("L65").Value = 24090.908
$ws.Range("M65").Value = -720827.8999999999
$ws.Range("N65").Value = -30330.908

$ws.Range("H107").Value = 38036.293
$ws.Range("I107").Value = 3343.1667
$ws.Range("K107").Value = 10029.5001
$ws.Range("M107").Value = -8109.500100000001

$ws.Range("H132").Value = 9546.066000000001
$ws.Range("I132").Value = 11154.728
$ws.Range("J132").Value = 5122.25
$ws.Range("K132").Value = 33464.18399999999
$ws.Range("L132").Value = 15366.75
$ws.Range("M132").Value = -30934.18399999999
$ws.Range("N132").Value = -20426.75

$ws.Range("H136").Value = 2203.2917
$ws.Range("I136").Value = 1543.0588
$ws.Range("K136").Value = 4629.1764
$ws.Range("M136").Value = -2079.1764
